# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Overview sheet: status text changes from "Ready for handoff" to
#    "Handed back: in sync with en-US" for both zh-cn and de-de rows.
#  - zh-cn / de-de sheets: the "Latest Target File" / "Latest Handback File"
#    / "Latest Handback DateTime" columns get populated with the real
#    handback file names + timestamps (previously empty / placeholder).
#  - New hyperlinks are added on the "Latest Target File" cells.
#  - A few columns are widened to fit the new, longer text.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet (status column for each language)
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

# Widen the zh-cn / de-de status columns to fit the longer text.
$overview.Range("E1").ColumnWidth = 29.14
$overview.Range("F1").ColumnWidth = 29.14

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText

$zhcn.Range("I2").Value = "a.md"
$zhcn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-26 10:35:35"

$zhcn.Range("I3").Value = "a.md"
$zhcn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-26 10:35:35"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e1fa9b6cdb870ce4f76ff8639e3bfd02fa518675/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md")
$zhcn.Range("I2").Style = "Hyperlink"

$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e1fa9b6cdb870ce4f76ff8639e3bfd02fa518675/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md")
$zhcn.Range("I3").Style = "Hyperlink"

$zhcn.Range("C1").ColumnWidth = 29.14
$zhcn.Range("J1").ColumnWidth = 39.15

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText

$dede.Range("I2").Value = "a.md"
$dede.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K2").Value = "2016-08-26 10:35:41"

$dede.Range("I3").Value = "a.md"
$dede.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K3").Value = "2016-08-26 10:35:41"

$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e1fa9b6cdb870ce4f76ff8639e3bfd02fa518675/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md")
$dede.Range("I2").Style = "Hyperlink"

$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e1fa9b6cdb870ce4f76ff8639e3bfd02fa518675/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md")
$dede.Range("I3").Style = "Hyperlink"

$dede.Range("C1").ColumnWidth = 29.14
$dede.Range("J1").ColumnWidth = 39.15

$wb.Save()
